$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in Unit_or_Levels column (C) and fix the id description (B2) ---
# --- and append the new education3 variable row (41) ---
$ws.Range("B2").Value = "Patient identifier from FMP "
$ws.Range("C2").Value = "Unique ID"
$ws.Range("C3").Value = "Female, Male, or Other"
$ws.Range("C4").Value = "American Indian and Alaska Native, Asian, Black or African American, More Than One Race, Native Hawaiian and Other Pacific Islander,  Other, Patient Refused, Unknown, White or Caucasian"
$ws.Range("C5").Value = "1= yes, 0 = no"
$ws.Range("C6").Value = "1=yes, 0 = no"
$ws.Range("C7").Value = "1= yes, 0 = no"
$ws.Range("C8").Value = "Count"
$ws.Range("C9").Value = "Numeric (years)"
$ws.Range("C11").Value = "No, Yes"
$ws.Range("C12").Value = "Medicaid/CHF, Medicare, Other, Private, Veteran/ASM"
$ws.Range("C13").Value = "Employed = Employed or Active Duty, Unemployed= Unemployed, Not in labor force= Student, Retired."
$ws.Range("C14").Value = "1= yes, 0 = no"
$ws.Range("C15").Value = "Has_Partner, Single"
$ws.Range("C16").Value = "< High School, High School Diploma, Some College, Partial College or 2 Year Degree, 4 Year College Degree, Graduate or Professional Degree, Unknown"
$ws.Range("C17").Value = "No, Yes"
$ws.Range("C18").Value = "Numeric (months)"
$ws.Range("C19").Value = "Numeric"
$ws.Range("C20").Value = "Numeric"
$ws.Range("C21").Value = "No, Yes"
$ws.Range("C22").Value = "Numeric, 0-10 possible values"
$ws.Range("C23").Value = "Numeric, 0-10 possible values"
$ws.Range("C24").Value = "Numeric, 0-10 possible values"
$ws.Range("C25").Value = "Numeric, 0-10 possible values"
$ws.Range("C26").Value = "Numeric, 0-10 possible values"
$ws.Range("C27").Value = "Numeric, 0-10 possible values"
$ws.Range("C28").Value = "Numeric, 0-10 possible values"
$ws.Range("C29").Value = "Numeric, 0-10 possible values"
$ws.Range("C30").Value = "Numeric, 0-80 possible values"
$ws.Range("C31").Value = "Numeric"
$ws.Range("C32").Value = "Numeric"
$ws.Range("C33").Value = "Numeric"
$ws.Range("C34").Value = "Numeric"
$ws.Range("C35").Value = "Other, White"
$ws.Range("C36").Value = "Disability, No Disability"
$ws.Range("C37").Value = "Numeric (monthly)"
$ws.Range("C38").Value = "Proportion, 0-1"
$ws.Range("A41").Value = "education3"
$ws.Range("B41").Value = "Time spent in education- Classification derived from Loucks, E.B., Abrahamowicz, M., Xiao, Y. et al. Associations of education with 30 year life course blood pressure trajectories: Framingham Offspring Study. BMC Public Health 11, 139 (2011). https://doi.org/10.1186/1471-2458-11-139"
$ws.Range("C41").Value = "≤12 years (reflecting high school or less), 13-16 years (some post-secondary education, including technical school and college degree), ≥17 years education (more than an undergraduate college degree)"

# --- Widen column C to fit the new Unit_or_Levels text ---
$ws.Columns.Item(3).ColumnWidth = 50.75

# --- Row heights to fit wrapped text in the new/updated cells ---
$ws.Rows.Item(3).RowHeight = 29.25
$ws.Rows.Item(4).RowHeight = 159
$ws.Rows.Item(5).RowHeight = 29.25
$ws.Rows.Item(11).RowHeight = 29.25
$ws.Rows.Item(12).RowHeight = 57.75
$ws.Rows.Item(13).RowHeight = 87
$ws.Rows.Item(16).RowHeight = 130.5
$ws.Rows.Item(22).RowHeight = 29.25
$ws.Rows.Item(23).RowHeight = 29.25
$ws.Rows.Item(24).RowHeight = 29.25
$ws.Rows.Item(25).RowHeight = 29.25
$ws.Rows.Item(26).RowHeight = 29.25
$ws.Rows.Item(27).RowHeight = 29.25
$ws.Rows.Item(28).RowHeight = 29.25
$ws.Rows.Item(29).RowHeight = 29.25
$ws.Rows.Item(30).RowHeight = 29.25
$ws.Rows.Item(36).RowHeight = 43.5
$ws.Rows.Item(37).RowHeight = 29.25
$ws.Rows.Item(39).RowHeight = 29.25
$ws.Rows.Item(40).RowHeight = 57.75
$ws.Rows.Item(41).RowHeight = 159

# --- Restore selection to C3, matching the saved view ---
$ws.Range("C3").Select() | Out-Null
